$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "boxen",
    "feuern",
    "platzen",
    "drücken",
    "sterben",
    "hoffen",
    "mögen",
    "stoßen",
    "ärgern",
    "heulen",
    "sperren",
    "spüren",
    "scheinen",
    "schließen",
    "garen",
    "drucken",
    "lockern",
    "fangen",
    "achten",
    "landen",
    "töten",
    "tauschen",
    "decken",
    "kriegen",
    "klettern",
    "wüten",
    "liefern",
    "feiern",
    "wählen",
    "opfern",
    "zielen",
    "spinnen",
    "sammeln",
    "zünden",
    "fallen",
    "wirken",
    "schwören",
    "knarren",
    "streichen",
    "treffen",
    "stecken",
    "warnen",
    "suchen",
    "fordern",
    "liegen",
    "wachsen",
    "flehen",
    "graben",
    "brauen",
    "brauchen",
    "gleichen",
    "dienen",
    "fischen",
    "gelten",
    "heilen",
    "rechnen",
    "tanzen",
    "fahren",
    "helfen",
    "morden",
    "malen",
    "nutzen",
    "bieten",
    "hören",
    "räumen",
    "sichern",
    "altern",
    "strahlen",
    "loben",
    "schätzen",
    "steuern",
    "wundern",
    "wiegen",
    "kichern",
    "dringen",
    "herzen",
    "tollen",
    "bitten",
    "schlucken",
    "werfen",
    "starren",
    "siegen",
    "drehen",
    "bergen",
    "beten",
    "klingen",
    "wehen",
    "tropfen",
    "trennen",
    "bauen",
    "stehlen",
    "flüchten",
    "zeigen",
    "nähen",
    "ehren",
    "schreiten",
    "scheitern",
    "machen",
    "schwächen",
    "pfeifen",
    "stürzen",
    "geben",
    "irren",
    "teilen",
    "reizen",
    "zögern",
    "fesseln",
    "lügen",
    "hupen",
    "trauen",
    "schenken",
    "regnen",
    "kranken",
    "bluten",
    "jagen",
    "lesen",
    "sorgen",
    "jubeln",
    "hassen",
    "spielen",
    "greifen",
    "arten",
    "treiben",
    "rauchen",
    "betteln",
    "rufen",
    "schrecken",
    "rühren",
    "enden",
    "wärmen",
    "wagen",
    "kümmern",
    "gnaden",
    "freuen",
    "weichen",
    "reiten",
    "schulden",
    "knurren",
    "seufzen",
    "trotzen",
    "bremsen",
    "gründen",
    "sprengen",
    "flüstern",
    "quälen",
    "sinken",
    "wellen",
    "stammen",
    "kehren",
    "öffnen",
    "runden",
    "erben",
    "wehtun",
    "stechen",
    "hauen",
    "schicken",
    "heben",
    "treten",
    "biegen",
    "formen",
    "binden",
    "folgen",
    "orten",
    "haben",
    "kosten",
    "ziehen",
    "holen",
    "filmen",
    "merken",
    "schwingen",
    "werden",
    "baden",
    "grüßen",
    "fließen",
    "ändern",
    "pflanzen",
    "saufen",
    "planen",
    "sparen",
    "mauern",
    "wenden",
    "kämpfen",
    "backen",
    "münzen",
    "bellen",
    "retten",
    "führen",
    "schmecken",
    "fällen",
    "rasen",
    "äußern",
    "stimmen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

Write-Output "Updated $($words.Length) word values in column B"
